# 07 - Error Handling.pptx : minor improvements to a couple of slides
#
# Slide 10 ("Method recover()"): split the sentence describing recover()'s
# behaviour so that a clarifying phrase "whose symbol is " is inserted in
# the middle of it.
#
# Slide 21 ("Only three methods..."): the first two runs of the opening
# sentence are merged back into a single run (no visible text change).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 - Content Placeholder 2 - paragraph 1
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(1)

# The whole sentence currently lives in run 3 of the paragraph:
#   "Method " + "recover()" + " implements error recovery by skipping
#   tokens until it finds one in the follow set of the nonterminal
#   defined by the rule."
$run10 = $para10.Runs(3)

# Shrink the existing run down to the part that stays before the new
# clarifying phrase ...
$run10.Text = " implements error recovery by skipping tokens until it finds one "

# ... then add the new phrase right after it ...
$inserted10 = $run10.InsertAfter("whose symbol is ")

# ... followed by the remainder of the original sentence.
$inserted10.InsertAfter("in the follow set of the nonterminal defined by the rule.") | Out-Null

# ---------------------------------------------------------------------
# Slide 21 - Content Placeholder 2 - paragraph 1
# ---------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$sh21 = $s21.Shapes.Item(2)
$tr21 = $sh21.TextFrame.TextRange
$para21 = $tr21.Paragraphs(1)

# Originally two separate runs: "Only three " and "methods throw a ".
# Remove the leading "Only three " text and re-add it immediately in
# front of "methods throw a ", which merges the wording back into a
# single run.
$full21 = $tr21.Text
$leadStart = $full21.IndexOf("Only three")
$leadText = "Only three "
$leadRange = $tr21.Characters($leadStart + 1, $leadText.Length)
$leadRange.Text = ""

$full21b = $tr21.Text
$tailStart = $full21b.IndexOf("methods throw a")
$tailText = "methods throw a "
$tailRange = $tr21.Characters($tailStart + 1, $tailText.Length)
$tailRange.InsertBefore($leadText) | Out-Null
